# Update "想去人数" (interest count) values in column F for both the
# "展览" and "全部类型" worksheets, as captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row -> New Value updates for worksheet "展览" (sheet1)
$exhibitionUpdates = @{
    4  = 2023
    5  = 324
    6  = 574
    7  = 92
    8  = 2055
    9  = 10456
    12 = 272
    13 = 199
    15 = 7335
    17 = 695
    18 = 172
    20 = 3281
}

# Row -> New Value updates for worksheet "全部类型" (sheet4)
$allTypesUpdates = @{
    4  = 2023
    5  = 324
    6  = 574
    8  = 92
    9  = 2055
    12 = 10456
    15 = 272
    16 = 199
    18 = 7335
    20 = 695
    21 = 172
    23 = 3281
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
